# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new value for column F
$exhibitionUpdates = @{
    2  = 1167
    3  = 17
    4  = 1368
    5  = 309
    7  = 10655
    8  = 21
    10 = 289
    11 = 1035
    12 = 699
    13 = 12032
    14 = 12469
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": row -> new value for column F
$allTypesUpdates = @{
    3  = 1167
    4  = 17
    5  = 1368
    6  = 309
    8  = 10655
    9  = 21
    11 = 289
    12 = 1035
    13 = 699
    14 = 12032
    15 = 12469
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
